# Generate Report for Handback
# Update the timestamp strings recorded for the "c756cc0b-..." handback
# row (row 3) across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-03 16:52:00"

# --- zh-cn sheet ---
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-03 16:51:55"
$wsZhCn.Range("K3").Value = "2016-09-03 16:52:18"

# --- de-de sheet ---
# Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-09-03 16:52:25"
